$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the two transaction rows that remain
$ws.Range("A15").Value = "test"
$ws.Range("A16").Value = "test2"

# Remove the last two transaction rows (17 and 18) entirely
$ws.Rows("17:18").Delete()

# Leave the final selection on cell F1 (no frozen/scrolled top-left cell)
$ws.Activate()
$ws.Range("F1").Select()
